$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell reference -> new text value (values must remain literal text, matching
# the workbook author using plain-text cells for numeric-looking data).
$updates = @{
    'D2' = '311.90'
    'E2' = '8.03%'
    'G2' = '7'
    'D3' = '32.09'
    'E3' = '9.37%'
    'G3' = '7'
    'D4' = '5.337'
    'E4' = '5.14%'
    'G4' = '7'
    'D5' = '0.07561'
    'E5' = '13.10%'
    'G5' = '7'
    'D6' = '7.850'
    'E6' = '6.98%'
    'G6' = '7'
    'D7' = '3.689'
    'E7' = '8.33%'
    'G7' = '7'
    'D8' = '1.580'
    'E8' = '14.86%'
    'G8' = '7'
    'D9' = '0.9141'
    'E9' = '-0.26%'
    'G9' = '7'
    'D10' = '0.01698'
    'E10' = '2,529.63%'
    'G10' = '7'
    'D11' = '0.1704'
    'E11' = '7.22%'
    'G11' = '7'
    'D12' = '0.07695'
    'E12' = '12.83%'
    'G12' = '7'
    'D13' = '0.08152'
    'E13' = '6.43%'
    'G13' = '7'
    'D14' = '0.03023'
    'E14' = '2.88%'
    'G14' = '7'
    'D15' = '0.09879'
    'E15' = '10.00%'
    'G15' = '7'
    'D16' = '0.001525'
    'E16' = '-3.65%'
    'G16' = '7'
    'D17' = '0.04565'
    'E17' = '1.39%'
    'G17' = '7'
    'D18' = '0.006565'
    'E18' = '4.21%'
    'G18' = '7'
    'D19' = '3.504'
    'E19' = '1.50%'
    'G19' = '7'
    'D20' = '2.239'
    'E20' = '0.79%'
    'G20' = '7'
    'D21' = '0.3267'
    'E21' = '1.69%'
    'G21' = '7'
    'D22' = '0.1316'
    'E22' = '0.56%'
    'G22' = '7'
    'D23' = '4.177'
    'E23' = '2.78%'
    'G23' = '7'
    'E24' = '2.86%'
    'G24' = '7'
    'D25' = '0.001218'
    'E25' = '2.39%'
    'G25' = '7'
    'D26' = '0.004492'
    'E26' = '9.02%'
    'G26' = '7'
    'D27' = '0.0001299'
    'E27' = '8.33%'
    'G27' = '7'
    'E28' = '7.47%'
    'G28' = '7'
    'G29' = '7'
    'G30' = '7'
    'G31' = '7'
    'G32' = '7'
    'G33' = '7'
    'G34' = '7'
    'G35' = '7'
    'G36' = '7'
    'G37' = '7'
    'G38' = '7'
    'G39' = '7'
    'D40' = '0.04597'
    'E40' = '8.52%'
    'G40' = '7'
    'D41' = '0.007238'
    'E41' = '7.72%'
    'G41' = '7'
    'D42' = '0.1368'
    'E42' = '10.32%'
    'G42' = '7'
    'D43' = '0.002258'
    'E43' = '8.05%'
    'G43' = '7'
    'D44' = '0.01409'
    'E44' = '5.15%'
    'G44' = '7'
    'D45' = '0.00006149'
    'E45' = '7.38%'
    'G45' = '7'
    'G46' = '7'
    'E47' = '-0.65%'
    'G47' = '7'
    'G48' = '7'
    'G49' = '7'
    'G50' = '7'
    'G51' = '7'
}

foreach ($ref in $updates.Keys) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $updates[$ref]
    $rng.Style = "Normal"
}
